$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '65.759.12'
$ws.Range("E2").Value = '  -1.11%  '
$ws.Range("D3").Value = '3.456.89'
$ws.Range("E3").Value = '  -3.76%  '
$ws.Range("E4").Value = '  +0.00%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '596.65'
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = '  -1.71%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '137.44'
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = '  -7.52%  '
$ws.Range("D7").Value = '3.456.88'
$ws.Range("E7").Value = '  -3.72%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.494'
$ws.Range("D9").ClearFormats()
$ws.Range("E9").Value = '  +0.32%  '
$ws.Range("E10").Value = '  -5.99%  '
$ws.Range("E11").Value = '  -9.86%  '
$ws.Range("E12").Value = '  -7.80%  '
$ws.Range("D13").Value = '4.040.79'
$ws.Range("E13").Value = '  -3.80%  '
$ws.Range("E14").Value = '  -10.54%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '26.70'
$ws.Range("D15").ClearFormats()
$ws.Range("E15").Value = '  -9.89%  '
$ws.Range("B16").Value = 'WrappedEther'
$ws.Range("C16").Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range("D16").Value = '3.452.88'
$ws.Range("E16").Value = '  -3.38%  '
$ws.Range("B17").Value = 'WrappedBTC'
$ws.Range("C17").Value = 'https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc'
$ws.Range("D17").Value = '65.691.74'
$ws.Range("E17").Value = '  -1.30%  '
$ws.Range("E18").Value = '  -2.22%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '9.93'
$ws.Range("D19").ClearFormats()
$ws.Range("E19").Value = '  -10.34%  '
$ws.Range("E20").Value = '  -8.66%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '13.78'
$ws.Range("D21").ClearFormats()
$ws.Range("E21").Value = '  -7.34%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '396.14'
$ws.Range("D22").ClearFormats()
$ws.Range("E22").Value = '  -6.70%  '
$ws.Range("E23").Value = '  -10.31%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '73.58'
$ws.Range("D24").ClearFormats()
$ws.Range("E24").Value = '  -5.91%  '
$ws.Range("E25").Value = '  -0.01%  '
$ws.Range("D26").Value = '3.601.99'
$ws.Range("E27").Value = '  -10.66%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '7.23'
$ws.Range("D30").ClearFormats()
$ws.Range("E30").Value = '  -10.74%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '8.29'
$ws.Range("D31").ClearFormats()
$ws.Range("E31").Value = '  -11.64%  '
$ws.Range("D32").Value = '3.459.98'
$ws.Range("E32").Value = '  -3.65%  '
$ws.Range("E33").Value = '  +0.00%  '
$ws.Range("E34").Value = '  -6.69%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '22.99'
$ws.Range("D35").ClearFormats()
$ws.Range("E35").Value = '  -8.01%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '173.07'
$ws.Range("D36").ClearFormats()
$ws.Range("E36").Value = '  -1.17%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '1.23'
$ws.Range("D37").ClearFormats()
$ws.Range("E37").Value = '  -13.57%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '6.94'
$ws.Range("D38").ClearFormats()
$ws.Range("E38").Value = '  -10.48%  '
$ws.Range("E39").Value = '  -7.40%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '4.84'
$ws.Range("D40").ClearFormats()
$ws.Range("E40").Value = '  -12.56%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.0781'
$ws.Range("D41").ClearFormats()
$ws.Range("E41").Value = '  -8.49%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.824'
$ws.Range("D42").ClearFormats()
$ws.Range("E42").Value = '  -6.36%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '43.60'
$ws.Range("D43").ClearFormats()
$ws.Range("E43").Value = '  -5.38%  '
$ws.Range("E44").Value = '  +0.03%  '
$ws.Range("E45").Value = '  -14.28%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '1.64'
$ws.Range("D46").ClearFormats()
$ws.Range("E46").Value = '  -11.15%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '23.35'
$ws.Range("D47").ClearFormats()
$ws.Range("E47").Value = '  -2.64%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '1.11'
$ws.Range("D48").ClearFormats()
$ws.Range("E48").Value = '  -2.85%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '6.59'
$ws.Range("D49").ClearFormats()
$ws.Range("E49").Value = '  -7.68%  '
$ws.Range("E50").Value = '  -15.83%  '
$ws.Range("D51").Value = '2.222.54'
$ws.Range("E51").Value = '  -7.91%  '
